# Bold the resume's "headline" paragraphs: the name/title line, the
# section heading "Berufserfahrung", each job-title line, and the
# education section's degree line. Each of these paragraphs consists of
# a single run whose bold flag is explicitly turned off
# (<w:b w:val="0"/>); we flip it to bold (<w:b/>), matching the commit's
# intent of visually emphasizing those lines.

$d = $word.ActiveDocument

# Exact paragraph texts (as returned by Range.Text, which includes the
# trailing paragraph mark) that must become bold. Each appears once
# EXCEPT "Senior Animation Designer", which shows up twice as a whole
# paragraph: once as a standalone heading in the summary section
# (unchanged) and once as the actual job-title heading above the
# Contoso Animation entry from 2008-2014 (the one to bold). We
# disambiguate by walking the document's paragraphs in order and only
# bolding the requested (1st/2nd/...) whole-paragraph match.

function Set-ParagraphBold($paragraphRange) {
    $paragraphRange.Bold = 1
}

$targets = @(
    @{ Text = "Nestor Wilke"; Occurrence = 1 },
    @{ Text = "Berufserfahrung"; Occurrence = 1 },
    @{ Text = "Animation Team Manager"; Occurrence = 1 },
    @{ Text = "Senior Animation Designer"; Occurrence = 2 },
    @{ Text = "Animation Designer"; Occurrence = 1 },
    @{ Text = "Bachelor of Fine Arts in Animation"; Occurrence = 1 }
)

foreach ($target in $targets) {
    $seen = 0
    foreach ($p in $d.Paragraphs) {
        $text = $p.Range.Text
        # Paragraph text includes trailing paragraph-mark characters;
        # trim them (and surrounding whitespace) before comparing so we
        # match the whole-paragraph heading exactly (not a substring
        # inside a longer sentence).
        $trimmed = $text.TrimEnd([char]13, [char]7).Trim()
        if ($trimmed -eq $target.Text) {
            $seen = $seen + 1
            if ($seen -eq $target.Occurrence) {
                Set-ParagraphBold $p.Range
                break
            }
        }
    }
}
